# Lithuania A Lyga - daily odds-base refresh (11-04-2024 00:31)
#
# 1) Rows 26 and 27 (match ids 24 and 25) had their data swapped - the
#    "id" column (A) stays put, but every other field (match id, teams,
#    odds, etc.) moves from row 27 into row 26 and vice versa.
# 2) The last fixture row (row 136, id 134) was removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap the contents of row 26 and row 27 (columns B:AC) ---------

$row26 = 26
$row27 = 27
$firstCol = 2   # column B
$lastCol  = 29  # column AC

for ($col = $firstCol; $col -le $lastCol; $col++) {
    $cell26 = $ws.Cells.Item($row26, $col)
    $cell27 = $ws.Cells.Item($row27, $col)

    $val26 = $cell26.Value()
    $val27 = $cell27.Value()

    $cell26.Value = $val27
    $cell27.Value = $val26
}

# --- 2) Delete the final row (row 136) ---------------------------------

$ws.Rows.Item(136).Delete()
